$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D column values that look numeric stay as text, matching original inlineStr formatting
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '30.184.93'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '1.857.67'
$ws.Range('E3').Value = '  -2.84%  '
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '233.76'
$ws.Range('E5').Value = '  -2.64%  '
$ws.Range('D6').Value = '0.9994'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = '0.4696'
$ws.Range('E7').Value = '  -1.60%  '
$ws.Range('D8').Value = '0.2811'
$ws.Range('E8').Value = '  -1.03%  '
$ws.Range('D9').Value = '0.06547'
$ws.Range('E9').Value = '  -2.29%  '
$ws.Range('D10').Value = '20.04'
$ws.Range('E10').Value = '  +3.49%  '
$ws.Range('D11').Value = '0.07796'
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').Value = '97.04'
$ws.Range('E12').Value = '  -5.84%  '
$ws.Range('D13').Value = '1.861.99'
$ws.Range('E13').Value = '  -2.68%  '
$ws.Range('D14').Value = '5.099'
$ws.Range('E14').Value = '  -2.08%  '
$ws.Range('D15').Value = '0.6650'
$ws.Range('E15').Value = '  -0.92%  '
$ws.Range('D16').Value = '283.44'
$ws.Range('E16').Value = '  -2.12%  '
$ws.Range('D17').Value = '30.216.73'
$ws.Range('E17').Value = '  -1.02%  '
$ws.Range('D18').Value = '0.9990'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').Value = '5.451'
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').Value = '2.105.80'
$ws.Range('E21').Value = '  -2.59%  '
$ws.Range('D22').Value = '0.000007238'
$ws.Range('E22').Value = '  -3.23%  '
$ws.Range('D23').Value = '0.9995'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').Value = '6.141'
$ws.Range('E24').Value = '  -2.71%  '
$ws.Range('D25').Value = '167.43'
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('D26').Value = '9.303'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('D27').Value = '19.00'
$ws.Range('E27').Value = '  -2.03%  '
$ws.Range('E28').Value = '  -8.20%  '
$ws.Range('D29').Value = '1.339'
$ws.Range('E29').Value = '  -3.36%  '
$ws.Range('D30').Value = '0.09587'
$ws.Range('E30').Value = '  -3.50%  '
$ws.Range('D31').Value = '4.413'
$ws.Range('E31').Value = '  -3.27%  '
$ws.Range('D32').Value = '1.468'
$ws.Range('E32').Value = '  -3.13%  '
$ws.Range('D33').Value = '4.096'
$ws.Range('E33').Value = '  -3.70%  '
$ws.Range('D34').Value = '0.04672'
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('D35').Value = '1.099'
$ws.Range('E35').Value = '  -0.92%  '
$ws.Range('D36').Value = '0.6976'
$ws.Range('E36').Value = '  -3.74%  '
$ws.Range('D37').Value = '0.9984'
$ws.Range('D38').Value = '2.710'
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('E39').Value = '  -2.92%  '
$ws.Range('D40').Value = '6.327'
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('D41').Value = '2.507'
$ws.Range('E41').Value = '  -4.27%  '
$ws.Range('D42').Value = '71.97'
$ws.Range('E42').Value = '  -3.16%  '
$ws.Range('D43').Value = '0.8579'
$ws.Range('E43').Value = '  -0.98%  '
$ws.Range('D44').Value = '1.938'
$ws.Range('E44').Value = '  -0.69%  '
$ws.Range('D45').Value = '104.05'
$ws.Range('E45').Value = '  -1.62%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Value = '0.9990'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('B47').Value = 'TheSandbox'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D47').Value = '0.4158'
$ws.Range('E47').Value = '  -2.21%  '
$ws.Range('D48').Value = '1.012.92'
$ws.Range('E48').Value = '  +4.94%  '
$ws.Range('D49').Value = '7.194'
$ws.Range('E49').Value = '  -2.75%  '
$ws.Range('D50').Value = '9.000'
$ws.Range('E50').Value = '  +2.64%  '
$ws.Range('D51').Value = '33.69'
$ws.Range('E51').Value = '  -2.57%  '
